$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '70.887.27'
$ws.Cells.Item(2, 5).Value = '  +2.48%  '
$ws.Cells.Item(3, 4).Value = '3.808.24'
$ws.Cells.Item(3, 5).Value = '  +0.91%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '702.25'
$ws.Cells.Item(5, 5).Value = '  +10.91%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '172.91'
$ws.Cells.Item(6, 5).Value = '  +3.88%  '
$ws.Cells.Item(7, 4).Value = '3.806.42'
$ws.Cells.Item(7, 5).Value = '  +0.92%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 5).Value = '  +0.81%  '
$ws.Cells.Item(10, 5).Value = '  +2.15%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '7.67'
$ws.Cells.Item(11, 5).Value = '  +12.94%  '
$ws.Cells.Item(12, 5).Value = '  +0.34%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000254'
$ws.Cells.Item(13, 5).Value = '  +3.86%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '36.26'
$ws.Cells.Item(14, 5).Value = '  +3.74%  '
$ws.Cells.Item(15, 4).Value = '4.448.99'
$ws.Cells.Item(15, 5).Value = '  +0.97%  '
$ws.Cells.Item(16, 4).Value = '3.824.63'
$ws.Cells.Item(16, 5).Value = '  +2.07%  '
$ws.Cells.Item(17, 4).Value = '70.853.81'
$ws.Cells.Item(17, 5).Value = '  +2.43%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '17.71'
$ws.Cells.Item(18, 5).Value = '  +0.37%  '
$ws.Cells.Item(19, 5).Value = '  +2.65%  '
$ws.Cells.Item(20, 5).Value = '  +0.31%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '11.38'
$ws.Cells.Item(21, 5).Value = '  +19.11%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '480.22'
$ws.Cells.Item(22, 5).Value = '  +3.54%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.717'
$ws.Cells.Item(23, 5).Value = '  +1.16%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '83.89'
$ws.Cells.Item(24, 5).Value = '  +1.34%  '
$ws.Cells.Item(25, 5).Value = '  +0.12%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '12.36'
$ws.Cells.Item(26, 5).Value = '  +2.08%  '
$ws.Cells.Item(27, 5).Value = '  +1.00%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.47'
$ws.Cells.Item(28, 5).Value = '  +3.51%  '
$ws.Cells.Item(29, 4).Value = '3.959.13'
$ws.Cells.Item(29, 5).Value = '  +0.92%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '3.11'
$ws.Cells.Item(31, 5).Value = '  +15.70%  '
$ws.Cells.Item(32, 5).Value = '  +1.54%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '7.49'
$ws.Cells.Item(33, 5).Value = '  +5.43%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '29.56'
$ws.Cells.Item(34, 5).Value = '  +3.64%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.177'
$ws.Cells.Item(35, 5).Value = '  +5.35%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '9.21'
$ws.Cells.Item(36, 5).Value = '  +2.43%  '
$ws.Cells.Item(37, 5).Value = '  +0.09%  '
$ws.Cells.Item(38, 4).Value = '3.757.91'
$ws.Cells.Item(38, 5).Value = '  +0.86%  '
$ws.Cells.Item(39, 5).Value = '  +1.40%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.51'
$ws.Cells.Item(40, 5).Value = '  +5.26%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '5.99'
$ws.Cells.Item(41, 5).Value = '  +3.06%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.000336'
$ws.Cells.Item(42, 5).Value = '  +23.80%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.21'
$ws.Cells.Item(43, 5).Value = '  +12.42%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.970'
$ws.Cells.Item(44, 5).Value = '  +0.75%  '
$ws.Cells.Item(45, 5).Value = '  +0.02%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '45.51'
$ws.Cells.Item(47, 5).Value = '  +5.22%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '160.50'
$ws.Cells.Item(48, 5).Value = '  +1.39%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '48.76'
$ws.Cells.Item(49, 5).Value = '  +4.30%  '
$ws.Cells.Item(50, 5).Value = '  -1.16%  '
$ws.Cells.Item(51, 5).Value = '  +1.05%  '

# Restore default (Normal) style for cells we temporarily formatted as Text,
# so no residual style index diverges from the original workbook formatting.
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(49, 4).Style = 'Normal'
